$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "scenarioA" to "ScenarioA".
# Excel automatically keeps the defined names that reference the sheet
# (_xlnm._FilterDatabase, businfo) in sync with the new name.
$ws.Name = "ScenarioA"

# Update the header labels in row 3: the short codes "lat"/"lon" are
# replaced with the descriptive column names "Latitude"/"Longitude".
$ws.Range("M3").Value = "Latitude"
$ws.Range("N3").Value = "Longitude"
